$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# Update the last-updated timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 22:15"

# Direct numeric updates for individual country rows (Espana, Suiza, Noruega, Japon)
Set-RowValues 7 @(18077, 3308, 1107, 16139, 939, 193, 831)
Set-RowValues 12 @(4133, 1018, 15, 4075, 0, 10, 43)
Set-RowValues 17 @(1790, 199, 1, 1782, 27, 1, 7)
Set-RowValues 20 @(943, 29, 191, 719, 46, 4, 33)

# Turquia climbs the ranking with updated case data; rows 34-52 (Polonia..Rusia)
# cascade down by one position and row 53 now shows India (was row 52)
$ws.Cells.Item(34, 1).Value = "Turquia"
Set-RowValues 34 @(359, 168, 0, 355, 0, 2, 4)
$ws.Cells.Item(35, 1).Value = "Polonia"
Set-RowValues 35 @(355, 68, 13, 337, 3, 0, 5)
$ws.Cells.Item(36, 1).Value = "Singapur"
Set-RowValues 36 @(345, 32, 124, 221, 14, 0, 0)
$ws.Cells.Item(37, 1).Value = "Chile"
Set-RowValues 37 @(342, 104, 0, 342, 6, 0, 0)
$ws.Cells.Item(38, 1).Value = "Luxemburgo"
Set-RowValues 38 @(335, 132, 6, 325, 1, 2, 4)
$ws.Cells.Item(39, 1).Value = "Islandia"
Set-RowValues 39 @(330, 80, 5, 325, 1, 0, 0)
$ws.Cells.Item(40, 1).Value = "Eslovenia"
Set-RowValues 40 @(319, 33, 0, 318, 6, 0, 1)
$ws.Cells.Item(41, 1).Value = "Indonesia"
Set-RowValues 41 @(308, 81, 15, 268, 0, 6, 25)
$ws.Cells.Item(42, 1).Value = "Barein"
Set-RowValues 42 @(278, 22, 110, 167, 4, 0, 1)
$ws.Cells.Item(43, 1).Value = "Rumania"
Set-RowValues 43 @(277, 17, 25, 252, 5, 0, 0)
$ws.Cells.Item(44, 1).Value = "Arabia Saudita"
Set-RowValues 44 @(274, 36, 8, 266, 0, 0, 0)
$ws.Cells.Item(45, 1).Value = "Tailandia"
Set-RowValues 45 @(272, 60, 42, 229, 1, 0, 1)
$ws.Cells.Item(46, 1).Value = "Estonia"
Set-RowValues 46 @(267, 9, 1, 266, 1, 0, 0)
$ws.Cells.Item(47, 1).Value = "Egipto"
Set-RowValues 47 @(256, 46, 42, 207, 0, 1, 7)
$ws.Cells.Item(48, 1).Value = "Peru"
Set-RowValues 48 @(234, 89, 1, 233, 7, 0, 0)
$ws.Cells.Item(49, 1).Value = "Filipinas"
Set-RowValues 49 @(217, 15, 8, 192, 1, 0, 17)
$ws.Cells.Item(50, 1).Value = "Hong Kong"
Set-RowValues 50 @(208, 15, 98, 106, 4, 0, 4)
$ws.Cells.Item(51, 1).Value = "Ecuador"
Set-RowValues 51 @(199, 31, 1, 195, 2, 0, 3)
$ws.Cells.Item(52, 1).Value = "Rusia"
Set-RowValues 52 @(199, 52, 8, 190, 0, 1, 1)
$ws.Cells.Item(53, 1).Value = "India"
Set-RowValues 53 @(194, 25, 20, 170, 0, 1, 4)

# Republica de Macedonia (row 87) updated figures
Set-RowValues 87 @(50, 7, 1, 49, 1, 0, 0)
